# Update average_county_temperature (column AA) for rows 9-36 with the
# new NOAA-sourced temperature value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 19.30324074074072

for ($row = 9; $row -le 36; $row++) {
    $ws.Cells.Item($row, 27).Value = $newValue
}
